{"js": "// Update the date line and the 25 two-digit multiplication problems to the\n// values from the next day's generated worksheet. Each old value is unique\n// in the document, so a scoped search-and-replace per pair is unambiguous.\n// Replacements are applied in document (top-to-bottom) order so that a\n// newly-written value (e.g. the \"56\u00d766=\" produced in row 4) can never be\n// mistaken for a not-yet-processed original cell later in the same pass.\nconst replacements = [\n  [\"2025-02-09 Sunday\", \"2025-02-10 Monday\"],\n  [\"82\u00d744=\", \"47\u00d732=\"],\n  [\"54\u00d790=\", \"56\u00d777=\"],\n  [\"69\u00d768=\", \"29\u00d756=\"],\n  [\"63\u00d755=\", \"85\u00d744=\"],\n  [\"71\u00d779=\", \"11\u00d743=\"],\n  [\"90\u00d757=\", \"60\u00d788=\"],\n  [\"56\u00d766=\", \"71\u00d758=\"],\n  [\"37\u00d718=\", \"23\u00d791=\"],\n  [\"78\u00d729=\", \"55\u00d724=\"],\n  [\"99\u00d753=\", \"91\u00d775=\"],\n  [\"91\u00d769=\", \"16\u00d769=\"],\n  [\"71\u00d739=\", \"31\u00d724=\"],\n  [\"86\u00d749=\", \"71\u00d781=\"],\n  [\"70\u00d758=\", \"15\u00d785=\"],\n  [\"75\u00d776=\", \"56\u00d766=\"],\n  [\"31\u00d787=\", \"44\u00d798=\"],\n  [\"51\u00d747=\", \"27\u00d776=\"],\n  [\"45\u00d767=\", \"84\u00d769=\"],\n  [\"19\u00d764=\", \"72\u00d731=\"],\n  [\"49\u00d760=\", \"93\u00d784=\"],\n  [\"76\u00d714=\", \"42\u00d724=\"],\n  [\"35\u00d739=\", \"12\u00d778=\"],\n  [\"60\u00d797=\", \"94\u00d713=\"],\n  [\"28\u00d775=\", \"25\u00d731=\"],\n  [\"25\u00d739=\", \"63\u00d719=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  // The source text is unique per cell, so only the first hit is touched;\n  // this avoids ever renaming a cell twice within one pass.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 two-digit multiplication problems to the\n# values from the next day's generated worksheet. Each old value is unique\n# in the document, so Find.Execute (which stops at the first match) can be\n# used as a safe scoped replace for each pair.\n# Replacements are applied in document (top-to-bottom) order so that a\n# newly-written value (e.g. the \"56\u00d766=\" produced in row 4) can never be\n# mistaken for a not-yet-processed original cell later in the same pass.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-09 Sunday\", \"2025-02-10 Monday\"),\n    @(\"82\u00d744=\", \"47\u00d732=\"),\n    @(\"54\u00d790=\", \"56\u00d777=\"),\n    @(\"69\u00d768=\", \"29\u00d756=\"),\n    @(\"63\u00d755=\", \"85\u00d744=\"),\n    @(\"71\u00d779=\", \"11\u00d743=\"),\n    @(\"90\u00d757=\", \"60\u00d788=\"),\n    @(\"56\u00d766=\", \"71\u00d758=\"),\n    @(\"37\u00d718=\", \"23\u00d791=\"),\n    @(\"78\u00d729=\", \"55\u00d724=\"),\n    @(\"99\u00d753=\", \"91\u00d775=\"),\n    @(\"91\u00d769=\", \"16\u00d769=\"),\n    @(\"71\u00d739=\", \"31\u00d724=\"),\n    @(\"86\u00d749=\", \"71\u00d781=\"),\n    @(\"70\u00d758=\", \"15\u00d785=\"),\n    @(\"75\u00d776=\", \"56\u00d766=\"),\n    @(\"31\u00d787=\", \"44\u00d798=\"),\n    @(\"51\u00d747=\", \"27\u00d776=\"),\n    @(\"45\u00d767=\", \"84\u00d769=\"),\n    @(\"19\u00d764=\", \"72\u00d731=\"),\n    @(\"49\u00d760=\", \"93\u00d784=\"),\n    @(\"76\u00d714=\", \"42\u00d724=\"),\n    @(\"35\u00d739=\", \"12\u00d778=\"),\n    @(\"60\u00d797=\", \"94\u00d713=\"),\n    @(\"28\u00d775=\", \"25\u00d731=\"),\n    @(\"25\u00d739=\", \"63\u00d719=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
